$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New task names typed into column B for rows 15-23
$ws.Range("B15").Value = "Create behavior for highlighted active tab of modal"
$ws.Range("B16").Value = "Creating classes and interfaces for tabs of modal"
$ws.Range("B17").Value = "Applying Thymeleaf for Data Profile"
$ws.Range("B18").Value = "Test add row if working"
$ws.Range("B19").Value = "Re-import csv file to postgres"
$ws.Range("B20").Value = "Research about Event Handling JS"
$ws.Range("B21").Value = "Apply Event Handling Method in JS for ID retrieve"
$ws.Range("B22").Value = "Apply Ajax call for ID retrieve"
$ws.Range("B23").Value = "Remodify code for Controller for Model And View"

# New logged hours across the daily columns
$ws.Range("R14").Value = 5

$ws.Range("R15").Value = 3
$ws.Range("S15").Value = 5

$ws.Range("S16").Value = 3

$ws.Range("T17").Value = 8
$ws.Range("U17").Value = 6

$ws.Range("X18").Value = 2

$ws.Range("X19").Value = 6
$ws.Range("Y19").Value = 4

$ws.Range("Y20").Value = 2

$ws.Range("Y21").Value = 2
$ws.Range("Z21").Value = 8

$ws.Range("AA22").Value = 8
$ws.Range("AB22").Value = 2

$ws.Range("AB23").Value = 6

$ws.Range("Q39").Value = 1
$ws.Range("S39").Value = 1

# Move the cursor/selection to where the author last left it
$null = $ws.Range("K22").Select()
